$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial (45180 -> 2023-09-11) for every
# data row (rows 2 through 199). Bump it by one day (45181 -> 2023-09-12)
# for all of them, matching the commit's "Automatic update of files."
$ws.Range("C2:C199").Value = 45181
